$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TC6_SearchResults_Typeahead")
$ws2 = $wb.Worksheets.Item("Testdata")

$nbsp = [char]0x00A0
$ws2.Range("B6").Value = "Showing Results for" + $nbsp + "`"sprocket`""

$ws1.Range("B11").Select()
$ws2.Activate()
$ws2.Range("B6").Select()
